# Update gh-pages to output generated at 456a3b4
# Updates "想去人数" (want-to-go count) figures in the 展览 and 全部类型 sheets.

$wb = $excel.ActiveWorkbook

# --- Sheet "展览" (index 1) ---
$wsExpo = $wb.Worksheets.Item("展览")
$wsExpo.Range("F2").Value = 3384
$wsExpo.Range("F3").Value = 19
$wsExpo.Range("F4").Value = 66
$wsExpo.Range("F5").Value = 1556
$wsExpo.Range("F6").Value = 54
$wsExpo.Range("F7").Value = 326

# --- Sheet "全部类型" (index 4) ---
$wsAll = $wb.Worksheets.Item("全部类型")
$wsAll.Range("F2").Value = 3384
$wsAll.Range("F3").Value = 19
$wsAll.Range("F4").Value = 66
$wsAll.Range("F5").Value = 1556
$wsAll.Range("F6").Value = 54
$wsAll.Range("F8").Value = 326
